$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
# Row 86 (G context=12603)
$ws.Range("H86").Value = 1831.9166
$ws.Range("I86").Value = 1419.4
$ws.Range("K86").Value = 1419.4
$ws.Range("M86").Value = -296.4000000000001
# Row 89 (G context=12603)
$ws.Range("H89").Value = 1831.9166
$ws.Range("I89").Value = 1419.4
$ws.Range("K89").Value = 7097
$ws.Range("M89").Value = -1481
# Row 134 (G context=41997)
$ws.Range("H134").Value = 45000
$ws.Range("J134").Value = 45000
$ws.Range("L134").Value = 45000
$ws.Range("N134").Value = -55140
# Row 135 (G context=44047)
$ws.Range("H135").Value = 2668.516
$ws.Range("J135").Value = 11187.5
$ws.Range("L135").Value = 100687.5
$ws.Range("N135").Value = -105757.5

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
# Row 61 (G context=43999)
$ws.Range("H61").Value = 3084.8076
$ws.Range("I61").Value = 3446.647
$ws.Range("J61").Value = 2401.3333
$ws.Range("K61").Value = 3446.647
$ws.Range("L61").Value = 2401.3333
$ws.Range("M61").Value = -3234.647
$ws.Range("N61").Value = -2825.3333
# Row 74 (G context=44000)
$ws.Range("H74").Value = 1265.7778
$ws.Range("I74").Value = 778.4
$ws.Range("J74").Value = 1875
$ws.Range("K74").Value = 778.4
$ws.Range("L74").Value = 1875
$ws.Range("M74").Value = 95.60000000000002
$ws.Range("N74").Value = -3623
# Row 77 (G context=44000)
$ws.Range("H77").Value = 1265.7778
$ws.Range("I77").Value = 778.4
$ws.Range("J77").Value = 1875
$ws.Range("K77").Value = 3892
$ws.Range("L77").Value = 9375
$ws.Range("M77").Value = 476
$ws.Range("N77").Value = -18111
# Row 102 (G context=19945)
$ws.Range("H102").Value = 4992
$ws.Range("I102").Value = 5478.8887
$ws.Range("J102").Value = 4115.6
$ws.Range("K102").Value = 5478.8887
$ws.Range("L102").Value = 4115.6
$ws.Range("M102").Value = -3856.8887
$ws.Range("N102").Value = -7359.6
# Row 119 (G context=26287)
$ws.Range("H119").Value = 31749.5
$ws.Range("J119").Value = 31749.5
$ws.Range("L119").Value = 31749.5
$ws.Range("N119").Value = -41425.5
# Row 132 (G context=43997)
$ws.Range("H132").Value = 1923.566
$ws.Range("I132").Value = 1786.7297
$ws.Range("J132").Value = 2240
$ws.Range("K132").Value = 5360.189100000001
$ws.Range("L132").Value = 6720
$ws.Range("M132").Value = -2830.189100000001
$ws.Range("N132").Value = -11780
# Row 136 (G context=43999)
$ws.Range("H136").Value = 3084.8076
$ws.Range("I136").Value = 3446.647
$ws.Range("J136").Value = 2401.3333
$ws.Range("K136").Value = 10339.941
$ws.Range("L136").Value = 7203.999899999999
$ws.Range("M136").Value = -7789.940999999999
$ws.Range("N136").Value = -12303.9999

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
# Row 29 (G context=2408)
$ws.Range("H29").Value = 10000
$ws.Range("J29").Value = 10000
$ws.Range("L29").Value = 10000
$ws.Range("N29").Value = -10586
# Row 31 (G context=44023)
$ws.Range("H31").Value = 4372.616
$ws.Range("I31").Value = 2748.795
$ws.Range("J31").Value = 5720.0425
$ws.Range("K31").Value = 2748.795
$ws.Range("L31").Value = 5720.0425
$ws.Range("M31").Value = -2453.795
$ws.Range("N31").Value = -6310.0425
# Row 34 (G context=44023)
$ws.Range("H34").Value = 4372.616
$ws.Range("I34").Value = 2748.795
$ws.Range("J34").Value = 5720.0425
$ws.Range("K34").Value = 2748.795
$ws.Range("L34").Value = 5720.0425
$ws.Range("M34").Value = -2546.795
$ws.Range("N34").Value = -6124.0425
# Row 50 (G context=1862)
$ws.Range("H50").Value = 8231.833000000001
$ws.Range("J50").Value = 8231.833000000001
$ws.Range("L50").Value = 8231.833000000001
$ws.Range("N50").Value = -9481.833000000001
# Row 51 (G context=2039)
$ws.Range("H51").Value = 8728.143
$ws.Range("I51").Value = 5000
$ws.Range("J51").Value = 9349.5
$ws.Range("K51").Value = 5000
$ws.Range("L51").Value = 9349.5
$ws.Range("M51").Value = -4264
$ws.Range("N51").Value = -10821.5
# Row 60 (G context=1937)
$ws.Range("H60").Value = 8042.2856
$ws.Range("J60").Value = 8200.5
$ws.Range("L60").Value = 8200.5
$ws.Range("N60").Value = -9222.5
# Row 61 (G context=2039)
$ws.Range("H61").Value = 8728.143
$ws.Range("I61").Value = 5000
$ws.Range("J61").Value = 9349.5
$ws.Range("K61").Value = 5000
$ws.Range("L61").Value = 9349.5
$ws.Range("M61").Value = -4652
$ws.Range("N61").Value = -10045.5
# Row 68 (G context=10611)
$ws.Range("H68").Value = 17100.334
$ws.Range("J68").Value = 17100.334
$ws.Range("L68").Value = 17100.334
$ws.Range("N68").Value = -18598.334
# Row 71 (G context=10611)
$ws.Range("H71").Value = 17100.334
$ws.Range("J71").Value = 17100.334
$ws.Range("L71").Value = 51301.00199999999
$ws.Range("N71").Value = -58789.00199999999
# Row 74 (G context=10636)
$ws.Range("H74").Value = 16344.5
$ws.Range("J74").Value = 16344.5
$ws.Range("L74").Value = 16344.5
$ws.Range("N74").Value = -18092.5
# Row 77 (G context=10636)
$ws.Range("H77").Value = 16344.5
$ws.Range("J77").Value = 16344.5
$ws.Range("L77").Value = 49033.5
$ws.Range("N77").Value = -57769.5
# Row 105 (G context=19928)
$ws.Range("H105").Value = 726.8
$ws.Range("I105").Value = 585.3333
$ws.Range("K105").Value = 585.3333
$ws.Range("M105").Value = 1161.6667
# Row 132 (G context=44019)
$ws.Range("H132").Value = 2268.0667
$ws.Range("I132").Value = 1421.55
$ws.Range("J132").Value = 3961.1
$ws.Range("K132").Value = 4264.65
$ws.Range("L132").Value = 11883.3
$ws.Range("M132").Value = -1734.65
$ws.Range("N132").Value = -16943.3
# Row 140 (G context=42455)
$ws.Range("H140").Value = 90000
$ws.Range("J140").Value = 90000
$ws.Range("L140").Value = 90000
$ws.Range("N140").Value = -100360
# Row 141 (G context=43345)
$ws.Range("H141").Value = 18236.8
$ws.Range("I141").Value = 10296
$ws.Range("J141").Value = 50000
$ws.Range("K141").Value = 10296
$ws.Range("L141").Value = 50000
$ws.Range("M141").Value = -5116
$ws.Range("N141").Value = -60360

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
# Row 132 (G context=44008)
$ws.Range("H132").Value = 2493.6897
$ws.Range("I132").Value = 1785.1111
$ws.Range("J132").Value = 3653.182
$ws.Range("K132").Value = 5355.3333
$ws.Range("L132").Value = 10959.546
$ws.Range("M132").Value = -2825.3333
$ws.Range("N132").Value = -16019.546

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
# Row 7 (G context=36249)
$ws.Range("H7").Value = 2364
$ws.Range("I7").Value = 1910.7097
$ws.Range("J7").Value = 3535
$ws.Range("K7").Value = 1910.7097
$ws.Range("L7").Value = 3535
$ws.Range("M7").Value = -1798.7097
$ws.Range("N7").Value = -3759
# Row 40 (G context=36248)
$ws.Range("H40").Value = 2673.8333
$ws.Range("I40").Value = 2383.3157
$ws.Range("J40").Value = 3777.8
$ws.Range("K40").Value = 2383.3157
$ws.Range("L40").Value = 3777.8
$ws.Range("M40").Value = -2247.3157
$ws.Range("N40").Value = -4049.8
# Row 119 (G context=26288)
$ws.Range("H119").Value = 29866.666
$ws.Range("J119").Value = 29866.666
$ws.Range("L119").Value = 29866.666
$ws.Range("N119").Value = -39542.666
# Row 126 (G context=36249)
$ws.Range("H126").Value = 2364
$ws.Range("I126").Value = 1910.7097
$ws.Range("J126").Value = 3535
$ws.Range("K126").Value = 5732.1291
$ws.Range("L126").Value = 10605
$ws.Range("M126").Value = -3262.1291
$ws.Range("N126").Value = -15545

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
# Row 54 (G context=3413)
$ws.Range("H54").Value = 2000
$ws.Range("J54").Value = 0
$ws.Range("L54").Value = 0
$ws.Range("N54").ClearContents()
# Row 119 (G context=26289)
$ws.Range("H119").Value = 33579.2
$ws.Range("J119").Value = 33579.2
$ws.Range("L119").Value = 33579.2
$ws.Range("N119").Value = -43255.2
# Row 132 (G context=44029)
$ws.Range("H132").Value = 23812454
$ws.Range("I132").Value = 31251030
$ws.Range("J132").Value = 9007.532999999999
$ws.Range("K132").Value = 93753090
$ws.Range("L132").Value = 27022.599
$ws.Range("M132").Value = -93750560
$ws.Range("N132").Value = -32082.599
